$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.527.74"
$ws.Range("E2").Value = "  +3.77%  "
$ws.Range("D3").Value = "1.739.92"
$ws.Range("E3").Value = "  +4.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.35%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4813"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2682"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06245"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("D10").Value = "1.738.38"
$ws.Range("E10").Value = "  +4.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07132"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6223"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.540"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9999"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "26.536.13"
$ws.Range("E17").Value = "  +3.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006896"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.65%  "
$ws.Range("D21").Value = "1.960.90"
$ws.Range("E21").Value = "  +4.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.592"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.902"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.364"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.815"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.419"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.006"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.747"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07881"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04591"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.19%  "
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.004"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6388"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9293"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "112.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.24%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.995"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.93%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.435"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.004"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01521"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.759"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3922"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.015"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1205"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05330"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.895"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.36%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.254"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3449"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.64%  "
